$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: remove "PU32" from D1, add "Fig in crore" to J1
$ws.Range("D1").Value = $null
$ws.Range("J1").Value = "Fig in crore"

# Row 3: new headers
$ws.Range("B3").Value = "2021-22"
$ws.Range("C3").Value = "2022-23"
$ws.Range("D3").Value = "JUN' 21"
$ws.Range("E3").Value = "JUN' 22"
$ws.Range("F3").Value = "JUN' 22"

# Row 4 (SMH 01)
$ws.Range("B4").Value = 7.72
$ws.Range("C4").Value = 6.29
$ws.Range("D4").Value = 1.58
$ws.Range("E4").Value = 1.64
$ws.Range("F4").Value = 2.7

# Row 5 (SMH 02)
$ws.Range("B5").Value = 90.73
$ws.Range("C5").Value = 79.2
$ws.Range("D5").Value = 21.18
$ws.Range("E5").Value = 20.59
$ws.Range("F5").Value = 29.56

# Row 6 (SMH 03)
$ws.Range("B6").Value = 11.38
$ws.Range("C6").Value = 9.039999999999999
$ws.Range("D6").Value = 0.7
$ws.Range("E6").Value = 2.35
$ws.Range("F6").Value = 3.6

# Row 7 (SMH 04)
$ws.Range("B7").Value = 12.51
$ws.Range("C7").Value = 11.67
$ws.Range("D7").Value = 2.79
$ws.Range("E7").Value = 3.03
$ws.Range("F7").Value = 3.02

# Row 8 (SMH 05)
$ws.Range("B8").Value = 75.31
$ws.Range("C8").Value = 68.11
$ws.Range("D8").Value = 20.45
$ws.Range("E8").Value = 17.71
$ws.Range("F8").Value = 18.74

# Row 9 (SMH 06)
$ws.Range("B9").Value = 68.48
$ws.Range("C9").Value = 53.89
$ws.Range("D9").Value = 13.66
$ws.Range("E9").Value = 14.01
$ws.Range("F9").Value = 17.39

# Row 10 (SMH 07)
$ws.Range("B10").Value = 37.83
$ws.Range("C10").Value = 30.43
$ws.Range("D10").Value = 7.8
$ws.Range("E10").Value = 7.91
$ws.Range("F10").Value = 12.12

# Row 11 (SMH 08)
$ws.Range("B11").Value = 0

# Row 12 (SMH 09)
$ws.Range("B12").Value = 37.84
$ws.Range("C12").Value = 32.36
$ws.Range("D12").Value = 8.890000000000001
$ws.Range("E12").Value = 8.41
$ws.Range("F12").Value = 11.31

# Row 13 (SMH 10)
$ws.Range("B13").Value = 8.98
$ws.Range("C13").Value = 7.87
$ws.Range("D13").Value = 2.63
$ws.Range("E13").Value = 2.05
$ws.Range("F13").Value = 4.33

# Row 14 (SMH 11)
$ws.Range("B14").Value = 0

# Row 15 (TOTAL)
$ws.Range("B15").Value = 350.79
$ws.Range("C15").Value = 298.86
$ws.Range("D15").Value = 79.69
$ws.Range("E15").Value = 77.7
$ws.Range("F15").Value = 102.77
